# Regenerate the handoff report:
#   - the "ready for handoff" source file's uuid was re-minted
#     (b589d2d7-e3e4-45f5-a0a5-d611781d48c4 -> 55002fb5-9332-468a-ac8e-f6f9be16b2cf)
#   - its handoff package hash / handoff timestamps were refreshed
#   - the old "Handoff transform failed" row (09738bc4-...) no longer
#     applies and drops out of the report, so the trailing
#     ".localization-config" row shifts up one row on every sheet.

$wb = $excel.ActiveWorkbook

$newGuid = "55002fb5-9332-468a-ac8e-f6f9be16b2cf"
$newHash = "f53d1ae240a912c97e30dc541e6dc2c2e57333ef"

$baseUrl = "https://github.com/OpenLocalizationTest/oltest/blob/bd0e827ff61a21e7235cacefff79562915686f8c"

# ---------------------------------------------------------------
# Sheet "Overview": File Name / zh-cn / de-de summary
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 3 ("09738bc4....md" / "Handoff transform failed") is gone; row 4
# (".localization-config" / "Not to be localized") slides up into row 3.
$wsOverview.Range("A3:C3").Value2 = $wsOverview.Range("A4:C4").Value2
$wsOverview.Rows.Item(4).Delete()

# Row 2 keeps its shape, just the source file's guid changed.
$wsOverview.Range("A2").Value2 = "$newGuid.md"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "$baseUrl/e2e/$newGuid.md", "", "", "$newGuid.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "$baseUrl/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------
# Sheet "zh-cn": per-language handoff detail
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A3").Value2 = $wsZh.Range("A4").Value2
$wsZh.Range("B3").Value2 = $wsZh.Range("B4").Value2
$wsZh.Range("C3").ClearContents()
$wsZh.Range("D3").Value2 = $wsZh.Range("D4").Value2
$wsZh.Range("G3").Value2 = $wsZh.Range("G4").Value2
$wsZh.Range("H3").Value2 = $wsZh.Range("H4").Value2
$wsZh.Rows.Item(4).Delete()

$wsZh.Range("A2").Value2 = "$newGuid.md"
$wsZh.Range("C2").Value2 = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("D2").Value2 = "2016-01-09 04:43:58"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "$baseUrl/e2e/$newGuid.md", "", "", "$newGuid.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1c7f00186822fa64d6ba92e41a73c0bf712c3137/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/$newGuid.$newHash.zh-cn.xlf", "", "", "$newGuid.$newHash.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "$baseUrl/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------
# Sheet "de-de": per-language handoff detail
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A3").Value2 = $wsDe.Range("A4").Value2
$wsDe.Range("B3").Value2 = $wsDe.Range("B4").Value2
$wsDe.Range("C3").ClearContents()
$wsDe.Range("D3").Value2 = $wsDe.Range("D4").Value2
$wsDe.Range("G3").Value2 = $wsDe.Range("G4").Value2
$wsDe.Range("H3").Value2 = $wsDe.Range("H4").Value2
$wsDe.Rows.Item(4).Delete()

$wsDe.Range("A2").Value2 = "$newGuid.md"
$wsDe.Range("C2").Value2 = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("D2").Value2 = "2016-01-09 04:44:07"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "$baseUrl/e2e/$newGuid.md", "", "", "$newGuid.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3d8b4514a7b3b3254176848f9171416cf67db7ac/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/$newGuid.$newHash.de-de.xlf", "", "", "$newGuid.$newHash.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "$baseUrl/.localization-config", "", "", ".localization-config")

Write-Output "Report regenerated for handoff"
